$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "thang 6" (June) ratio column J and cap marker column K for the relevant rows
$ws.Range("J2").Value = 0.7
$ws.Range("K2").Value = "*"

$ws.Range("J3").Value = 0.5
$ws.Range("K3").Value = "*"

$ws.Range("J4").Value = 0.7
$ws.Range("K4").Value = "*"

$ws.Range("J6").Value = 0.7
$ws.Range("K6").Value = "*"

$ws.Range("J7").Value = 0.7
$ws.Range("K7").Value = "*"

# Reflect the user's final selection on the active pane
$ws.Range("J3").Select()
